$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (from H1, which already has the bold/border/centered style)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(1,8,6,8,7,1,6,4,7,6,8,10,8,6,7,7,8,8,7,7,7,7,6,7,6,9,8,6,7,7,8,7,8,7,8,6,9,7,7,6,5,6,6,4,6,6,7,7,4,5,5,4,6,6,8,7,7,5,5,4,8,7,4,5)
$jValues = @(1,8,6,8,7,2,6,5,7,6,8,10,8,7,7,7,8,8,8,8,7,7,6,8,7,9,8,6,7,8,9,7,8,8,8,6,9,8,7,6,6,7,6,5,7,6,7,8,5,6,6,6,7,7,9,8,7,5,5,5,8,7,4,5)

for ($r = 2; $r -le 65; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
